$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update column F ("想去人数") values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3565
$ws1.Range("F5").Value = 3565
$ws1.Range("F6").Value = 259
$ws1.Range("F7").Value = 5084
$ws1.Range("F8").Value = 5084
$ws1.Range("F9").Value = 521
$ws1.Range("F10").Value = 353
$ws1.Range("F11").Value = 199
$ws1.Range("F12").Value = 691
$ws1.Range("F14").Value = 86
$ws1.Range("F16").Value = 696
$ws1.Range("F18").Value = 34
$ws1.Range("F23").Value = 4906
$ws1.Range("F24").Value = 4906
$ws1.Range("F28").Value = 6032
$ws1.Range("F32").Value = 338
$ws1.Range("F38").Value = 1007
$ws1.Range("F42").Value = 869
$ws1.Range("F43").Value = 979
$ws1.Range("F44").Value = 2027

# Sheet "全部类型" (fourth sheet) - update column F ("想去人数") values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1117
$ws4.Range("F7").Value = 3565
$ws4.Range("F8").Value = 3565
$ws4.Range("F9").Value = 259
$ws4.Range("F10").Value = 5084
$ws4.Range("F11").Value = 5084
$ws4.Range("F12").Value = 521
$ws4.Range("F13").Value = 353
$ws4.Range("F14").Value = 199
$ws4.Range("F15").Value = 691
$ws4.Range("F17").Value = 86
$ws4.Range("F19").Value = 696
$ws4.Range("F21").Value = 34
$ws4.Range("F27").Value = 4906
$ws4.Range("F28").Value = 4906
$ws4.Range("F32").Value = 6032
$ws4.Range("F36").Value = 338
$ws4.Range("F43").Value = 1007
$ws4.Range("F47").Value = 869
$ws4.Range("F48").Value = 979
$ws4.Range("F50").Value = 2027
